$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing unstyled data cell as the template style so that newly
# entered text-like values (e.g. "300.18", "-0.36%") keep the default style
# instead of picking up an automatic 'Text' number format style.
$baseStyle = $ws.Range("F2").Style

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = $baseStyle
}

Set-TextValue "D2" "300.18"
Set-TextValue "E2" "-0.36%"
Set-TextValue "G2" "21"
Set-TextValue "D3" "31.81"
Set-TextValue "E3" "1.37%"
Set-TextValue "G3" "21"
Set-TextValue "E4" "0.72%"
Set-TextValue "G4" "21"
Set-TextValue "D5" "0.08143"
Set-TextValue "E5" "10.52%"
Set-TextValue "G5" "21"
Set-TextValue "D6" "2.430"
Set-TextValue "E6" "1.30%"
Set-TextValue "G6" "21"
Set-TextValue "D7" "7.839"
Set-TextValue "E7" "-1.40%"
Set-TextValue "G7" "21"
Set-TextValue "D8" "3.874"
Set-TextValue "E8" "2.16%"
Set-TextValue "G8" "21"
Set-TextValue "D9" "0.9258"
Set-TextValue "E9" "0.92%"
Set-TextValue "G9" "21"
Set-TextValue "D10" "0.1760"
Set-TextValue "E10" "2.53%"
Set-TextValue "G10" "21"
Set-TextValue "D11" "0.07395"
Set-TextValue "E11" "-3.11%"
Set-TextValue "G11" "21"
Set-TextValue "D12" "0.09034"
Set-TextValue "E12" "11.51%"
Set-TextValue "G12" "21"
Set-TextValue "D13" "0.03030"
Set-TextValue "E13" "0.32%"
Set-TextValue "G13" "21"
Set-TextValue "D14" "0.1001"
Set-TextValue "E14" "0.82%"
Set-TextValue "G14" "21"
Set-TextValue "D15" "0.001520"
Set-TextValue "E15" "1.73%"
Set-TextValue "G15" "21"
Set-TextValue "D16" "0.005973"
Set-TextValue "E16" "-3.02%"
Set-TextValue "G16" "21"
Set-TextValue "D17" "3.591"
Set-TextValue "E17" "3.96%"
Set-TextValue "G17" "21"
Set-TextValue "D18" "2.286"
Set-TextValue "E18" "2.66%"
Set-TextValue "G18" "21"
Set-TextValue "D19" "0.3263"
Set-TextValue "E19" "-0.97%"
Set-TextValue "G19" "21"
Set-TextValue "D20" "0.1340"
Set-TextValue "E20" "0.27%"
Set-TextValue "G20" "21"
Set-TextValue "D21" "4.099"
Set-TextValue "E21" "-11.79%"
Set-TextValue "G21" "21"
Set-TextValue "D22" "0.1682"
Set-TextValue "E22" "7.44%"
Set-TextValue "G22" "21"
Set-TextValue "D23" "0.04630"
Set-TextValue "E23" "-0.27%"
Set-TextValue "G23" "21"
Set-TextValue "D24" "0.001248"
Set-TextValue "E24" "1.90%"
Set-TextValue "G24" "21"
Set-TextValue "D25" "0.004546"
Set-TextValue "E25" "1.36%"
Set-TextValue "G25" "21"
Set-TextValue "D26" "0.0001201"
Set-TextValue "E26" "-7.52%"
Set-TextValue "G26" "21"
Set-TextValue "D27" "0.0003411"
Set-TextValue "E27" "27.90%"
Set-TextValue "G27" "21"
Set-TextValue "G28" "21"
Set-TextValue "G29" "21"
Set-TextValue "G30" "21"
Set-TextValue "G31" "21"
Set-TextValue "G32" "21"
Set-TextValue "G33" "21"
Set-TextValue "G34" "21"
Set-TextValue "G35" "21"
Set-TextValue "G36" "21"
Set-TextValue "G37" "21"
Set-TextValue "G38" "21"
Set-TextValue "D39" "0.01762"
Set-TextValue "E39" "1.35%"
Set-TextValue "G39" "21"
Set-TextValue "D40" "0.04605"
Set-TextValue "E40" "1.80%"
Set-TextValue "G40" "21"
Set-TextValue "D41" "0.006911"
Set-TextValue "E41" "-4.38%"
Set-TextValue "G41" "21"
Set-TextValue "E42" "2.53%"
Set-TextValue "G42" "21"
Set-TextValue "D43" "0.002211"
Set-TextValue "E43" "-0.70%"
Set-TextValue "G43" "21"
Set-TextValue "D44" "0.009860"
Set-TextValue "E44" "-8.26%"
Set-TextValue "G44" "21"
Set-TextValue "D45" "0.00006310"
Set-TextValue "E45" "0.69%"
Set-TextValue "G45" "21"
Set-TextValue "E46" "0.02%"
Set-TextValue "G46" "21"
Set-TextValue "B47" "CoinbaseStockToken"
Set-TextValue "C47" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D47" "0.008403"
Set-TextValue "E47" "-15.90%"
Set-TextValue "G47" "21"
Set-TextValue "B48" "BOLO"
Set-TextValue "C48" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D48" "0.8025"
Set-TextValue "E48" "-0.77%"
Set-TextValue "G48" "21"
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.02%"
Set-TextValue "G49" "21"
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.09%"
Set-TextValue "G50" "21"
Set-TextValue "G51" "21"
